# Sync automatico del tracker: append new match rows (493-499) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{ Row = 493; EventId = "14494886"; Fecha = "2025-08-28"; JugadorA = "Jenson Brooksby"; JugadorB = "Flavio Cobolli";    Pronostico = "Gana Flavio Cobolli";     Cuota = 1.57 },
    @{ Row = 494; EventId = "14494990"; Fecha = "2025-08-28"; JugadorA = "Anna Kalinskaya";  JugadorB = "Yulia Putintseva";  Pronostico = "Gana Yulia Putintseva";    Cuota = 3.75 },
    @{ Row = 495; EventId = "14494986"; Fecha = "2025-08-28"; JugadorA = "Magdalena Frech";  JugadorB = "Peyton Stearns";    Pronostico = "Gana Magdalena Frech";     Cuota = 2.75 },
    @{ Row = 496; EventId = "14495030"; Fecha = "2025-08-28"; JugadorA = "Hailey Baptiste";  JugadorB = "Naomi Osaka";       Pronostico = "Gana Hailey Baptiste";     Cuota = 4 },
    @{ Row = 497; EventId = "14495031"; Fecha = "2025-08-28"; JugadorA = "Daria Kasatkina";  JugadorB = "Kamilla Rakhimova"; Pronostico = "Gana Kamilla Rakhimova";   Cuota = 2.75 },
    @{ Row = 498; EventId = "14487597"; Fecha = "2025-08-28"; JugadorA = "Stuart Parker";    JugadorB = "Jakub Paul";        Pronostico = "Gana Stuart Parker";       Cuota = 2.25 },
    @{ Row = 499; EventId = "14487554"; Fecha = "2025-08-29"; JugadorA = "Sanhui Shin";      JugadorB = "Kokoro Isomura";    Pronostico = "Gana Kokoro Isomura";      Cuota = 2.5 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # event_id / fecha look like numbers or dates, so force plain text
    # formatting before writing them, matching the text values already
    # used throughout column A/B in this tracker.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $r.EventId

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $r.Fecha

    $ws.Cells.Item($row, 3).Value = $r.JugadorA
    $ws.Cells.Item($row, 4).Value = $r.JugadorB
    $ws.Cells.Item($row, 5).Value = $r.Pronostico
    $ws.Cells.Item($row, 6).Value = $r.Cuota
}
